# RSTK-8093_RSTK-8109_RSTK-8939: add Stock Loc No column and new MS-project /
# Pro-1 lot-tracking sample rows to the "Location Add" template; keep
# Sheet2's layout (which already had the Stock Loc No column) as the
# reference shape.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 ("Location Add"): insert a new "Stock Loc No" column before
# the old column E (Quantity_Serial Nos), shifting E:T -> F:U. ---
$ws1.Columns("E:E").Insert()
$ws1.Range("E1").Value = "Stock Loc No"

# --- Insert two new rows above the old row 12 (multidiv serial / OH loc),
# shifting the old rows 12-21 down to 14-23. ---
$ws1.Rows("12:13").Insert()

# New row 12: multidiv serial (Lot and serial track), 1019 (100 MS-project),
# OH (On Hand Loc), Stock Loc No "Pro-1", Qty 10, Lot No "LT3"
$ws1.Range("A12").Value = "10 (Denver)"
$ws1.Range("B12").Value = "multidiv serial (Lot and serial track)"
$ws1.Range("C12").Value = "1019 (100 MS-project)"
$ws1.Range("D12").Value = "OH (On Hand Loc)"
$ws1.Range("E12").Value = "Pro-1"
$ws1.Range("F12").Value = 10
$ws1.Range("J12").Value = "LT3"

# New row 13: multidiv serial (Lot and serial track), 100 Home Project,
# OH (On Hand Loc), Stock Loc No "Pro-1", Qty 10, Lot No "LT1"
$ws1.Range("A13").Value = "10 (Denver)"
$ws1.Range("B13").Value = "multidiv serial (Lot and serial track)"
$ws1.Range("C13").Value = "100 Home Project"
$ws1.Range("D13").Value = "OH (On Hand Loc)"
$ws1.Range("E13").Value = "Pro-1"
$ws1.Range("F13").Value = 10
$ws1.Range("J13").Value = "LT1"

# --- View/selection state: Sheet2 gets a simple single-cell selection,
# Sheet1 gets a whole-row selection and remains the active (tab-selected)
# sheet. Sheet2's selection is set first so activating Sheet1 afterward
# leaves it as the active tab. ---
$ws2.Range("G15").Select()

$ws1.Activate()
$ws1.Range("A3:XFD3").Select()
